$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.198.58"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.852.24"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'0.6970"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").Value = "'237.64"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.07850"
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("D9").Value = "'0.3013"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("D10").Value = "'23.70"
$ws.Range("E10").Value = "  +2.15%  "
$ws.Range("D11").Value = "'0.08117"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "1.848.77"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").Value = "'5.175"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").Value = "'0.7039"
$ws.Range("E14").Value = "  -2.43%  "
$ws.Range("D15").Value = "'89.40"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "29.218.64"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "'5.794"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "'13.18"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").Value = "'235.60"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "2.101.25"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'7.502"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").Value = "'162.70"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("D26").Value = "'8.859"
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("D27").Value = "'0.1411"
$ws.Range("E27").Value = "  -0.89%  "
$ws.Range("D28").Value = "'17.99"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").Value = "'1.920"
$ws.Range("E29").Value = "  -1.73%  "
$ws.Range("D30").Value = "'1.403"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").Value = "'1.476"
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("D32").Value = "'4.311"
$ws.Range("E32").Value = "  -4.30%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").Value = "'0.05146"
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("D35").Value = "'1.166"
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("D36").Value = "'0.7081"
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("D37").Value = "'0.9982"
$ws.Range("E37").Value = "  -1.94%  "
$ws.Range("D38").Value = "'2.678"
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("D40").Value = "'2.703"
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("D41").Value = "1.151.27"
$ws.Range("E41").Value = "  +4.94%  "
$ws.Range("D42").Value = "'0.9227"
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("D43").Value = "'5.955"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "'0.4227"
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "'102.89"
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("D48").Value = "'0.5291"
$ws.Range("E48").Value = "  -2.94%  "
$ws.Range("E49").Value = "  -3.26%  "
$ws.Range("D50").Value = "'9.138"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").Value = "'6.938"
$ws.Range("E51").Value = "  -0.85%  "
